$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The "Price" column (D) holds values formatted like "27.287.92" which Excel
# would otherwise auto-coerce into a number (losing the textual formatting).
# Force the cells to Text format first so the literal string is preserved,
# matching the workbook's inlineStr (text) cells.
$priceCells = "D2","D3","D4","D5","D7","D8","D10","D11","D12","D13","D14","D15","D16","D17","D18","D19","D20","D21","D22","D23","D24","D25","D26","D27","D28","D29","D30","D31","D32","D33","D35","D36","D37","D38","D39","D41","D42","D43","D44","D45","D46","D48","D49","D50","D51"
foreach ($addr in $priceCells) {
    $ws.Range($addr).NumberFormat = "@"
}

# Row 2 - Bitcoin
$ws.Range("D2").Value = "27.287.92"
$ws.Range("E2").Value = "  +1.20%  "

# Row 3 - Ethereum
$ws.Range("D3").Value = "1.855.63"
$ws.Range("E3").Value = "  +1.62%  "

# Row 4 - TetherUSD
$ws.Range("D4").Value = "1.001"
$ws.Range("E4").Value = "  -0.57%  "

# Row 5 - BNB
$ws.Range("D5").Value = "314.02"
$ws.Range("E5").Value = "  +0.38%  "

# Row 6 - USDC
$ws.Range("E6").Value = "  -0.51%  "

# Row 7 - XRP
$ws.Range("D7").Value = "0.4609"
$ws.Range("E7").Value = "  +0.87%  "

# Row 8 - Cardano
$ws.Range("D8").Value = "0.3712"
$ws.Range("E8").Value = "  +0.38%  "

# Row 9 - Dogecoin
$ws.Range("E9").Value = "  -0.10%  "

# Row 10 - Polygon
$ws.Range("D10").Value = "0.8882"
$ws.Range("E10").Value = "  +1.86%  "

# Row 11 - Solana
$ws.Range("D11").Value = "20.18"
$ws.Range("E11").Value = "  +2.26%  "

# Row 12 - TRON
$ws.Range("D12").Value = "0.07838"
$ws.Range("E12").Value = "  -1.49%  "

# Row 13 - was Polkadot, now WrappedEther
$ws.Range("B13").Value = "WrappedEther"
$ws.Range("C13").Value = "https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth"
$ws.Range("D13").Value = "1.830.15"
$ws.Range("E13").Value = "  -6.80%  "

# Row 14 - was WrappedEther, now Polkadot
$ws.Range("B14").Value = "Polkadot"
$ws.Range("C14").Value = "https://coinranking.com/coin/25W7FG7om+polkadot-dot"
$ws.Range("D14").Value = "5.394"
$ws.Range("E14").Value = "  +1.32%  "

# Row 15 - Chainlink
$ws.Range("D15").Value = "6.538"
$ws.Range("E15").Value = "  -0.14%  "

# Row 16 - Litecoin
$ws.Range("D16").Value = "91.51"
$ws.Range("E16").Value = "  +0.19%  "

# Row 17 - BinanceUSD
$ws.Range("D17").Value = "1.002"
$ws.Range("E17").Value = "  -0.59%  "

# Row 18 - ShibaInu
$ws.Range("D18").Value = "0.000008934"
$ws.Range("E18").Value = "  +0.95%  "

# Row 19 - Dai
$ws.Range("D19").Value = "1.001"
$ws.Range("E19").Value = "  -0.57%  "

# Row 20 - Avalanche
$ws.Range("D20").Value = "14.78"
$ws.Range("E20").Value = "  +0.57%  "

# Row 21 - WrappedBTC
$ws.Range("D21").Value = "27.301.11"
$ws.Range("E21").Value = "  +0.71%  "

# Row 22 - Uniswap
$ws.Range("D22").Value = "5.115"
$ws.Range("E22").Value = "  +0.29%  "

# Row 23 - Cosmos
$ws.Range("D23").Value = "10.55"
$ws.Range("E23").Value = "  +0.22%  "

# Row 24 - WrappedliquidstakedEther2.0
$ws.Range("D24").Value = "2.058.57"
$ws.Range("E24").Value = "  -4.17%  "

# Row 25 - Toncoin
$ws.Range("D25").Value = "1.923"
$ws.Range("E25").Value = "  +4.21%  "

# Row 26 - Monero
$ws.Range("D26").Value = "151.97"
$ws.Range("E26").Value = "  -0.86%  "

# Row 27 - EthereumClassic
$ws.Range("D27").Value = "18.41"
$ws.Range("E27").Value = "  +0.26%  "

# Row 28 - LidoDAOToken
$ws.Range("D28").Value = "2.064"
$ws.Range("E28").Value = "  +1.31%  "

# Row 29 - BitcoinCash
$ws.Range("D29").Value = "115.96"
$ws.Range("E29").Value = "  +0.80%  "

# Row 30 - InternetComputer(DFINITY)
$ws.Range("D30").Value = "5.070"
$ws.Range("E30").Value = "  -1.34%  "

# Row 31 - Stellar
$ws.Range("D31").Value = "0.08822"
$ws.Range("E31").Value = "  -0.57%  "

# Row 32 - ImmutableX
$ws.Range("D32").Value = "0.7759"
$ws.Range("E32").Value = "  +6.73%  "

# Row 33 - HuobiToken
$ws.Range("D33").Value = "3.070"
$ws.Range("E33").Value = "  +3.54%  "

# Row 34 - ARBITRUM
$ws.Range("E34").Value = "  +3.85%  "

# Row 35 - was Filecoin, now RenderToken
$ws.Range("B35").Value = "RenderToken"
$ws.Range("C35").Value = "https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr"
$ws.Range("D35").Value = "2.773"
$ws.Range("E35").Value = "  +13.69%  "

# Row 36 - was RenderToken, now Filecoin
$ws.Range("B36").Value = "Filecoin"
$ws.Range("C36").Value = "https://coinranking.com/coin/ymQub4fuB+filecoin-fil"
$ws.Range("D36").Value = "4.500"
$ws.Range("E36").Value = "  +2.05%  "

# Row 37 - TrustWalletToken
$ws.Range("D37").Value = "1.083"
$ws.Range("E37").Value = "  +1.04%  "

# Row 38 - VeChain
$ws.Range("D38").Value = "0.01955"
$ws.Range("E38").Value = "  +0.93%  "

# Row 39 - Hedera
$ws.Range("D39").Value = "0.05251"
$ws.Range("E39").Value = "  +0.64%  "

# Row 40 - MXToken
$ws.Range("E40").Value = "  +0.48%  "

# Row 41 - FraxShare
$ws.Range("D41").Value = "7.066"
$ws.Range("E41").Value = "  -0.91%  "

# Row 42 - TheSandbox
$ws.Range("D42").Value = "0.5130"
$ws.Range("E42").Value = "  +0.05%  "

# Row 43 - Algorand
$ws.Range("D43").Value = "0.1638"
$ws.Range("E43").Value = "  +0.74%  "

# Row 44 - Aptos
$ws.Range("D44").Value = "8.419"
$ws.Range("E44").Value = "  +3.08%  "

# Row 45 - Decentraland
$ws.Range("D45").Value = "0.4808"
$ws.Range("E45").Value = "  -0.25%  "

# Row 46 - EnergySwap
$ws.Range("D46").Value = "10.39"
$ws.Range("E46").Value = "  +2.32%  "

# Row 47 - PaxDollar (unchanged)

# Row 48 - was Quant, now NEARProtocol
$ws.Range("B48").Value = "NEARProtocol"
$ws.Range("C48").Value = "https://coinranking.com/coin/DCrsaMv68+nearprotocol-near"
$ws.Range("D48").Value = "1.646"
$ws.Range("E48").Value = "  +1.07%  "

# Row 49 - was NEARProtocol, now Quant
$ws.Range("B49").Value = "Quant"
$ws.Range("C49").Value = "https://coinranking.com/coin/bauj_21eYVwso+quant-qnt"
$ws.Range("D49").Value = "102.27"
$ws.Range("E49").Value = "  -0.10%  "

# Row 50 - Cronos
$ws.Range("D50").Value = "0.06213"
$ws.Range("E50").Value = "  +0.08%  "

# Row 51 - Aave
$ws.Range("D51").Value = "65.70"
$ws.Range("E51").Value = "  +1.09%  "
